$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("AT2").Value = 2.62
$ws.Range("G4").Value = 1.8
$ws.Range("I4").Value = 4.75
$ws.Range("L4").Value = 5.5
$ws.Range("Z4").Value = 13
$ws.Range("AH4").Value = 10
$ws.Range("AI4").Value = 23
$ws.Range("AJ4").Value = 17
$ws.Range("AW4").Value = 6.5
$ws.Range("AX4").Value = 29
$ws.Range("G5").Value = 2.5
$ws.Range("H5").Value = 2.82
$ws.Range("K5").Value = 1.87
$ws.Range("Q5").Value = 2.5
$ws.Range("R5").Value = 1.5
$ws.Range("V5").Value = 1.67
$ws.Range("G6").Value = 2.15
$ws.Range("I6").Value = 3.7
$ws.Range("J6").Value = 2.88
$ws.Range("M6").Value = 1.08
$ws.Range("N6").Value = 8
$ws.Range("R6").Value = 1.54
$ws.Range("W6").Value = 6.5
$ws.Range("X6").Value = 9.5
$ws.Range("AA6").Value = 21
$ws.Range("AF6").Value = 51
$ws.Range("AH6").Value = 9
$ws.Range("AI6").Value = 17
$ws.Range("AJ6").Value = 13
$ws.Range("AU6").Value = 8.5
$ws.Range("AZ6").Value = 67
$ws.Range("M7").Value = 1.05
$ws.Range("O7").Value = 1.29
$ws.Range("M8").Value = 1.05
$ws.Range("O8").Value = 1.29
$ws.Range("O12").Value = 1.2
$ws.Range("P12").Value = 4.33
$ws.Range("Q12").Value = 1.67
$ws.Range("R12").Value = 2.15
$ws.Range("I13").Value = 1.67
$ws.Range("M13").Value = 1.01
$ws.Range("N13").Value = 23
$ws.Range("AB13").Value = 26
$ws.Range("AG13").Value = 81
$ws.Range("AQ13").Value = 51
$ws.Range("AX13").Value = 8.5
$ws.Range("AZ13").Value = 23
$ws.Range("G14").Value = 3.9
$ws.Range("I14").Value = 1.73
$ws.Range("Q14").Value = 1.36
$ws.Range("U14").Value = 1.4
$ws.Range("V14").Value = 2.75
$ws.Range("X14").Value = 26
$ws.Range("AB14").Value = 23
$ws.Range("AK14").Value = 17
$ws.Range("AO14").Value = 19
$ws.Range("BC14").Value = 201
$ws.Range("G15").Value = 1.72
$ws.Range("Q15").Value = 1.63
$ws.Range("G16").Value = 1.5
$ws.Range("N16").Value = 12
$ws.Range("Q16").Value = 1.77
$ws.Range("R16").Value = 1.97
$ws.Range("Q17").Value = 1.41
$ws.Range("O18").Value = 1.29
$ws.Range("P18").Value = 3.5
$ws.Range("Q18").Value = 1.87
$ws.Range("R18").Value = 1.87
$ws.Range("G19").Value = 1.72
$ws.Range("Q19").Value = 1.67
$ws.Range("L20").Value = 1.91
$ws.Range("J21").Value = 1.8
$ws.Range("K21").Value = 2.88
$ws.Range("N21").Value = 26
$ws.Range("Q21").Value = 1.33
$ws.Range("R21").Value = 3.4
$ws.Range("R24").Value = 1.57
$ws.Range("Q25").Value = 1.9
$ws.Range("R25").Value = 1.95
$ws.Range("Q26").Value = 1.77
$ws.Range("M27").Value = 1.08
$ws.Range("N27").Value = 8
$ws.Range("Q27").Value = 2.3
$ws.Range("J28").Value = 2.63
$ws.Range("M28").Value = 1.05
$ws.Range("N28").Value = 11
$ws.Range("Q28").Value = 1.98
$ws.Range("R28").Value = 1.88
$ws.Range("G30").Value = 2.05
$ws.Range("I30").Value = 3.6
$ws.Range("M30").Value = 1.05
$ws.Range("O30").Value = 1.27
$ws.Range("X30").Value = 10
$ws.Range("Z30").Value = 19
$ws.Range("AX30").Value = 19
$ws.Range("G32").Value = 2.38
$ws.Range("I32").Value = 2.7
$ws.Range("M32").Value = 1.02
$ws.Range("O32").Value = 1.13
$ws.Range("I33").Value = 2.3
$ws.Range("M33").Value = 1.05
$ws.Range("N33").Value = 8
$ws.Range("O33").Value = 1.41
$ws.Range("P33").Value = 2.62
$ws.Range("G34").Value = 1.62
$ws.Range("M34").Value = 1.03
$ws.Range("O34").Value = 1.22
$ws.Range("Q34").Value = 1.85
$ws.Range("R34").Value = 2
$ws.Range("Q36").Value = 1.98
$ws.Range("R36").Value = 1.83
$ws.Range("J39").Value = 2.88
